# Apply updated cryptocurrency price/volume data (and the row 47/48 swap)
# generated from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.930.74'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '3.765.65'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''640.96'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = '''165.50'
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = '3.767.21'
$ws.Range("E7").Value = '  -0.94%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '''0.524'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").Value = '''0.455'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '''6.92'
$ws.Range("E12").Value = '  +4.70%  '
$ws.Range("D13").Value = '''0.0000239'
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("D14").Value = '''34.91'
$ws.Range("E14").Value = '  -2.85%  '
$ws.Range("D15").Value = '4.398.51'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").Value = '3.772.00'
$ws.Range("E16").Value = '  +2.22%  '
$ws.Range("D17").Value = '68.900.38'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '''17.65'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").Value = '''472.47'
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").Value = '''9.57'
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("D23").Value = '''0.705'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -5.68%  '
$ws.Range("D25").Value = '''81.66'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").Value = '''12.15'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = '''10.06'
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '3.914.28'
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("D31").Value = '''2.68'
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("D32").Value = '''2.27'
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").Value = '''7.13'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").Value = '''28.54'
$ws.Range("E34").Value = '  -2.04%  '
$ws.Range("D35").Value = '''0.173'
$ws.Range("E35").Value = '  +16.20%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = '3.721.20'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = '''8.88'
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").Value = '  -5.31%  '
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '''0.955'
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D45").Value = '''45.16'
$ws.Range("E45").Value = '  +5.37%  '
$ws.Range("E46").Value = '  +4.30%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '''47.78'
$ws.Range("E47").Value = '  +2.03%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''155.47'
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("D49").Value = '''1.40'
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").Value = '''0.295'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").Value = '''8.37'
$ws.Range("E51").Value = '  -0.68%  '
